$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Update the Participant ID query (row 2, column B) with the new, fixed
# Neo4j Cypher query text (adds diagnosis/genomic_info traversal and sorts
# the collected sample ids).
$newParticipantQuery = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE s.study_name in ["GECCO OICR: Molecular Pathological Epidemiology of Colorectal Cancer"]
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id
LIMIT 100
"@

$ws.Range("B2").Value = $newParticipantQuery

# The longer query text needs more wrapped lines, so the row grows taller.
$ws.Rows.Item(2).RowHeight = 299.25

# Update the last active selection to C2, matching the author's saved state.
$ws.Range("C2").Select()

$wb.Save()
